$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.482.35"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.478.93"
$ws.Range("E3").Value = "  +9.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "2.853.77"
$ws.Range("E14").Value = "  +9.30%  "
$ws.Range("D15").Value = "2.484.19"
$ws.Range("E15").Value = "  +9.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.864"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.83%  "
$ws.Range("D18").Value = "46.386.97"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("E25").Value = "  +6.48%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.81%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.22%  "
$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +24.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0779"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.117"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0304"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.42%  "
$ws.Range("D43").Value = "2.010.28"
$ws.Range("E43").Value = "  +12.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +33.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.89%  "
$ws.Range("D50").Value = "2.718.57"
$ws.Range("E50").Value = "  +9.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.74%  "
